# Regenerate the "K" column (column G) of the save_data sheet with freshly
# simulated strikeout values (s_vals), replacing the old Strike# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> newly calculated K value (column G).
$kValues = [ordered]@{
    2  = 2
    4  = 0
    5  = 0
    6  = 2
    7  = 2
    8  = 1
    9  = 0
    10 = 2
    11 = 0
    12 = 2
    13 = 2
    14 = 0
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 0
    20 = 0
    21 = 2
    22 = 0
    23 = 1
    24 = 1
    25 = 1
    26 = 2
    27 = 1
    28 = 1
    29 = 0
    30 = 1
    31 = 0
    32 = 3
    33 = 0
    34 = 0
    35 = 2
    36 = 1
    38 = 2
    39 = 0
    40 = 1
    41 = 3
    42 = 0
    43 = 2
    44 = 0
    45 = 0
    46 = 1
    47 = 2
    48 = 2
    49 = 0
    50 = 0
    51 = 1
    52 = 0
    53 = 2
    54 = 0
    55 = 1
    56 = 1
    57 = 1
    58 = 1
    59 = 2
    60 = 0
    61 = 3
    62 = 1
    63 = 0
    64 = 2
    65 = 1
    66 = 3
    67 = 1
    68 = 1
    69 = 1
    70 = 1
    71 = 0
    72 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
